$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 69.95072642790045
$ws.Range("C2").Value = 76.85815826018292
$ws.Range("D2").Value = 69.28832658452248
$ws.Range("E2").Value = 83.236934797983

$ws.Range("B3").Value = 97.35782852822994
$ws.Range("C3").Value = 96.42719266622905
$ws.Range("D3").Value = 97.28406050719582
$ws.Range("E3").Value = 97.09517587248354

$ws.Range("B4").Value = 99.48166329888224
$ws.Range("C4").Value = 99.3411419777306
$ws.Range("D4").Value = 99.40534650347828
$ws.Range("E4").Value = 99.43709878099789

$ws.Range("B5").Value = 98.87575810940083
$ws.Range("C5").Value = 98.93598966755872
$ws.Range("D5").Value = 98.80888079405626
$ws.Range("E5").Value = 98.83680986153772

$ws.Range("B6").Value = 98.46062414827618
$ws.Range("C6").Value = 98.21069741918284
$ws.Range("D6").Value = 98.35153449997112
$ws.Range("E6").Value = 98.28375928631742

$ws.Range("B7").Value = 97.32911851896667
$ws.Range("C7").Value = 97.21617910747689
$ws.Range("D7").Value = 97.3672489486587
$ws.Range("E7").Value = 97.35479812538314

$ws.Range("B8").Value = 95.97488365675812
$ws.Range("C8").Value = 95.93806015431944
$ws.Range("D8").Value = 95.89558257907818
$ws.Range("E8").Value = 95.84201227139823
